$d = $word.ActiveDocument

# 1) Summary table: update TFS item description
$summaryTable = $d.Tables.Item(1)
$summaryTable.Cell(1, 2).Range.Find.Execute("25654 - sanitize data (input/output);", $true, $false, $false, $false, $false, $true, 1, $false, "25997 – Move email notification from UI to backend;", 2)

# 2) Log table: normalize existing last row (merge "TFS " + "25654..." into a single run)
$logTable = $d.Tables.Item(2)
$lastRowIndex = $logTable.Rows.Count
$logTable.Cell($lastRowIndex, 2).Range.Find.Execute("TFS 25654 - sanitize data (input/output);", $true, $false, $false, $false, $false, $true, 1, $false, "TFS 25654 - sanitize data (input/output);", 2)

# 3) Log table: append a new row documenting the 25997 change
$newRow = $logTable.Rows.Add()
$newRowIndex = $newRow.Index
$logTable.Cell($newRowIndex, 1).Range.Text = "01/09/2023"
$logTable.Cell($newRowIndex, 2).Range.Text = "TFS 25997 – Move email notification from UI to backend;"
$logTable.Cell($newRowIndex, 3).Range.Text = "Lili Huang"

# 4) Step table: bump changeset number from 52191 to 52392
$stepTable = $d.Tables.Item(3)
$stepTable.Cell(1, 2).Range.Find.Execute("191", $true, $false, $false, $false, $false, $true, 1, $false, "392", 2)
